# Refresh the live crypto price/volume snapshot on Sheet1 (columns B-E).
# Row 7 and row 8 also swap places (USDC <-> Solana) to reflect the new ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '37.919.59'
$ws.Range('E2').Value = '  +1.77%  '

# Row 3
$ws.Range('D3').Value = '2.103.19'
$ws.Range('E3').Value = '  +2.02%  '

# Row 4
$ws.Range('E4').Value = '  -0.02%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.47'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.44%  '

# Row 6
$ws.Range('E6').Value = '  +0.21%  '

# Row 7
$ws.Range('B7').Value = 'Solana'
$ws.Range('C7').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.15'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +1.75%  '

# Row 8
$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.04%  '

# Row 9
$ws.Range('E9').Value = '  +1.62%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0781'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +2.72%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.105'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +2.79%  '

# Row 12
$ws.Range('D12').Value = '2.400.31'
$ws.Range('E12').Value = '  +1.54%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.58'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.35%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.28'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +2.01%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.773'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.83%  '

# Row 16
$ws.Range('E16').Value = '  +1.82%  '

# Row 17
$ws.Range('D17').Value = '2.097.03'
$ws.Range('E17').Value = '  +1.80%  '

# Row 18
$ws.Range('D18').Value = '37.832.46'
$ws.Range('E18').Value = '  +1.66%  '

# Row 19
$ws.Range('E19').Value = '  -2.13%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '70.96'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +2.51%  '

# Row 21
$ws.Range('D21').Value = '0.0₃0826'
$ws.Range('E21').Value = '  +1.88%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '228.15'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.01%  '

# Row 23
$ws.Range('E23').Value = '  -0.04%  '

# Row 24
$ws.Range('E24').Value = '  -0.62%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.40'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.33%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '168.25'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.30%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.140'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +9.00%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.97'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +2.12%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.43'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.87%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.52'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +2.47%  '

# Row 31
$ws.Range('E31').Value = '  +1.30%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.65'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +4.04%  '

# Row 33
$ws.Range('E33').Value = '  +1.64%  '

# Row 34
$ws.Range('E34').Value = '  -0.02%  '

# Row 35
$ws.Range('E35').Value = '  +0.68%  '

# Row 36
$ws.Range('E36').Value = '  +5.35%  '

# Row 37
$ws.Range('E37').Value = '  +4.47%  '

# Row 38
$ws.Range('E38').Value = '  -0.23%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.42'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -4.69%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0993'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +6.56%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.95'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.18%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '97.66'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.23%  '

# Row 43
$ws.Range('E43').Value = '  +1.19%  '

# Row 44
$ws.Range('D44').Value = '1.458.91'
$ws.Range('E44').Value = '  -0.53%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.17'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.34%  '

# Row 46
$ws.Range('E46').Value = '  +4.25%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '15.77'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +4.27%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.07'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -6.19%  '

# Row 49
$ws.Range('E49').Value = '  +2.91%  '

# Row 50
$ws.Range('E50').Value = '  +2.24%  '

# Row 51
$ws.Range('D51').Value = '2.297.21'
$ws.Range('E51').Value = '  +2.03%  '
